$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows
# (rows 2-294) from 2023-09-06 (45175) to 2023-09-08 (45177).
$ws.Range("C2:C294").Value = 45177

# Row 294 picks up an explicit row height (matches the rest of the sheet).
$ws.Rows.Item(294).RowHeight = 15

# Append the new record row 295.
$ws.Cells.Item(295, 1).Value = "A 41434-2023"

$ws.Cells.Item(295, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(295, 2).Value = 45175

$ws.Cells.Item(295, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(295, 3).Value = 45177

$ws.Cells.Item(295, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(295, 5).Value = "VANSBRO"
# Column F (Markägare) is intentionally left blank for this row.
$ws.Cells.Item(295, 7).Value = 2.5
$ws.Cells.Item(295, 8).Value = 0
$ws.Cells.Item(295, 9).Value = 0
$ws.Cells.Item(295, 10).Value = 0
$ws.Cells.Item(295, 11).Value = 0
$ws.Cells.Item(295, 12).Value = 0
$ws.Cells.Item(295, 13).Value = 0
$ws.Cells.Item(295, 14).Value = 0
$ws.Cells.Item(295, 15).Value = 0
$ws.Cells.Item(295, 16).Value = 0
$ws.Cells.Item(295, 17).Value = 0

# Column R keeps the wrap-text style seen on every other data row, with no content.
$ws.Cells.Item(295, 18).WrapText = $true
